# Avances en creación de vistas, procedimientos y triggers
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Add the new table column "Encargado de la tarea" (extends table to A1:E9)
$newCol = $lo.ListColumns.Add()
$newCol.Range.Cells.Item(1,1).Value2 = "Encargado de la tarea"

# Fill in "Encargado de la tarea" column values
$ws.Range("E2").Value2 = "Facundo"
$ws.Range("E3").Value2 = "Leandro Correa"
$ws.Range("E4").Value2 = "Jesus"
$ws.Range("E5").Value2 = "Leandro Correa"

# Apply centered alignment to D1 and the whole "Encargado de la tarea" column (E1:E9),
# matching the header/data style used across the rest of the table (copy format from
# an already centered cell so we reuse the existing style instead of creating new ones)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null
$ws.Range("A2").Copy() | Out-Null
$ws.Range("E1:E9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Update the comment on "Mostrar contenido de una venta" (Vista 2) with join details
$ws.Range("C3").Value2 = "Debería mostrar todos los artículos que formaron parte de una venta`nJoin entre Ventas, ArtículoVenta y Artículo para obtener el nombre del artículo"
$ws.Range("C3").WrapText = $true

# Fill in Vista 3 description/comment (average of sales)
$ws.Range("B4").Value2 = "Mostrar el promedio de ventas realizadas"
$ws.Range("C4").Value2 = "Promedio de importe de ventas`nTabla Ventas con AVG de importeTotal"
$ws.Range("C4").WrapText = $true

# Extend the comment for Trigger 1 (AFTER UPDATE)
$ws.Range("C7").Value2 = "Sería un AFTER UPDATE, de la tabla Artículos"

# Column widths (closest achievable render widths to the target 14.43 / 20 char widths)
$ws.Columns.Item(4).ColumnWidth = 13.59
$ws.Columns.Item(5).ColumnWidth = 19.1

# Sheet view: zoom + selection
$ws.Application.ActiveWindow.Zoom = 130
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("C5").Select()
